# Auto-generated Excel COM-interop script
# Applies the numeric data refresh described by the commit:
#   "chore: update Sheets via scheduled runner"
# Updates computed market-price / profit columns (H, I, J, K, L, M, N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR worksheets.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
# Row 32: update H32, J32, L32, N32
$ws.Range("H32").Value = 1113
$ws.Range("J32").Value = 769.75
$ws.Range("L32").Value = 769.75
$ws.Range("N32").Value = -1421.75
# Row 115: update H115, J115, L115, N115
$ws.Range("H115").Value = 1473
$ws.Range("J115").Value = 666
$ws.Range("L115").Value = 1998
$ws.Range("N115").Value = -5132
# Row 132: update H132, I132, K132, M132
$ws.Range("H132").Value = 2924
$ws.Range("I132").Value = 2836.0188
$ws.Range("K132").Value = 8508.056399999999
$ws.Range("M132").Value = -5978.056399999999
# Row 138: update H138, I138, K138, M138
$ws.Range("H138").Value = 2413.394
$ws.Range("I138").Value = 1841.6072
$ws.Range("K138").Value = 5524.821599999999
$ws.Range("M138").Value = -384.8215999999993
# Row 141: update H141, I141, K141, M141
$ws.Range("H141").Value = 4052.6667
$ws.Range("I141").Value = 3823.3
$ws.Range("K141").Value = 11469.9
$ws.Range("M141").Value = -6289.900000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32: update H32, I32, K32, M32
$ws.Range("H32").Value = 7578061
$ws.Range("I32").Value = 3877114.8
$ws.Range("K32").Value = 3877114.8
$ws.Range("M32").Value = -3876827.8
# Row 45: update H45, I45, K45, M45
$ws.Range("H45").Value = 4453.8
$ws.Range("I45").Value = 3219.8572
$ws.Range("K45").Value = 3219.8572
$ws.Range("M45").Value = -2842.8572
# Row 61: update H61, I61, K61, M61
$ws.Range("H61").Value = 3647.611
$ws.Range("I61").Value = 3429
$ws.Range("K61").Value = 3429
$ws.Range("M61").Value = -3217
# Row 63: update H63, I63, K63, M63
$ws.Range("H63").Value = 3609
$ws.Range("I63").Value = 2483.1667
$ws.Range("K63").Value = 2483.1667
$ws.Range("M63").Value = -1797.1667
# Row 66: update H66, I66, K66, M66
$ws.Range("H66").Value = 3609
$ws.Range("I66").Value = 2483.1667
$ws.Range("K66").Value = 12415.8335
$ws.Range("M66").Value = -8983.833500000001
# Row 92: update H92, J92, L92, N92
$ws.Range("H92").Value = 66775
$ws.Range("J92").Value = 66775
$ws.Range("L92").Value = 66775
$ws.Range("N92").Value = -71767
# Row 132: update H132, I132, K132, M132
$ws.Range("H132").Value = 2440.8
$ws.Range("I132").Value = 1804.4706
$ws.Range("K132").Value = 5413.4118
$ws.Range("M132").Value = -2883.4118
# Row 136: update H136, I136, K136, M136
$ws.Range("H136").Value = 3647.611
$ws.Range("I136").Value = 3429
$ws.Range("K136").Value = 10287
$ws.Range("M136").Value = -7737

$ws = $wb.Worksheets.Item("BSM")
# Row 62: update H62, J62, L62, N62
$ws.Range("H62").Value = 128000
$ws.Range("J62").Value = 128000
$ws.Range("L62").Value = 128000
$ws.Range("N62").Value = -129372
# Row 65: update H65, J65, L65, N65
$ws.Range("H65").Value = 128000
$ws.Range("J65").Value = 128000
$ws.Range("L65").Value = 384000
$ws.Range("N65").Value = -390864
# Row 86: update H86, I86, K86, M86
$ws.Range("H86").Value = 2734.5454
$ws.Range("I86").Value = 3596.6667
$ws.Range("K86").Value = 3596.6667
$ws.Range("M86").Value = -2473.6667
# Row 89: update H89, I89, K89, M89
$ws.Range("H89").Value = 2734.5454
$ws.Range("I89").Value = 3596.6667
$ws.Range("K89").Value = 17983.3335
$ws.Range("M89").Value = -12367.3335
# Row 108: update H108, J108, L108, N108
$ws.Range("H108").Value = 55000
$ws.Range("J108").Value = 55000
$ws.Range("L108").Value = 55000
$ws.Range("N108").Value = -62680
# Row 134: update H134, I134, K134, M134
$ws.Range("H134").Value = 9201858
$ws.Range("I134").Value = 1833241.5
$ws.Range("K134").Value = 5499724.5
$ws.Range("M134").Value = -5497189.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16: update H16, I16, K16, M16
$ws.Range("H16").Value = 2441.625
$ws.Range("I16").Value = 2304.7144
$ws.Range("K16").Value = 2304.7144
$ws.Range("M16").Value = -2017.7144
# Row 31: update H31, I31, J31, K31, L31, M31, N31
$ws.Range("H31").Value = 3804.739
$ws.Range("I31").Value = 2324.8572
$ws.Range("J31").Value = 6106.778
$ws.Range("K31").Value = 2324.8572
$ws.Range("L31").Value = 6106.778
$ws.Range("M31").Value = -2029.8572
$ws.Range("N31").Value = -6696.778
# Row 34: update H34, I34, J34, K34, L34, M34, N34
$ws.Range("H34").Value = 3804.739
$ws.Range("I34").Value = 2324.8572
$ws.Range("J34").Value = 6106.778
$ws.Range("K34").Value = 2324.8572
$ws.Range("L34").Value = 6106.778
$ws.Range("M34").Value = -2122.8572
$ws.Range("N34").Value = -6510.778
# Row 58: update H58, I58, J58, K58, L58, M58, N58
$ws.Range("H58").Value = 2798.5334
$ws.Range("I58").Value = 2189.3
$ws.Range("J58").Value = 4017
$ws.Range("K58").Value = 2189.3
$ws.Range("L58").Value = 4017
$ws.Range("M58").Value = -1986.3
$ws.Range("N58").Value = -4423
# Row 113: update H113, I113, K113, M113
$ws.Range("H113").Value = 2441.625
$ws.Range("I113").Value = 2304.7144
$ws.Range("K113").Value = 2304.7144
$ws.Range("M113").Value = -134.7143999999998
# Row 122: update H122, I122, K122, M122
$ws.Range("H122").Value = 1391.6666
$ws.Range("I122").Value = 837.5
$ws.Range("K122").Value = 2512.5
$ws.Range("M122").Value = -62.5
# Row 132: update H132, I132, K132, M132
$ws.Range("H132").Value = 1815.4706
$ws.Range("I132").Value = 1724.2
$ws.Range("K132").Value = 5172.6
$ws.Range("M132").Value = -2642.6
# Row 134: update H134, I134, J134, K134, L134, M134, N134
$ws.Range("H134").Value = 3878.3845
$ws.Range("I134").Value = 2713.3333
$ws.Range("J134").Value = 6499.75
$ws.Range("K134").Value = 8139.999899999999
$ws.Range("L134").Value = 19499.25
$ws.Range("M134").Value = -5604.999899999999
$ws.Range("N134").Value = -24569.25
# Row 136: update H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 2798.5334
$ws.Range("I136").Value = 2189.3
$ws.Range("J136").Value = 4017
$ws.Range("K136").Value = 6567.900000000001
$ws.Range("L136").Value = 4017
$ws.Range("M136").Value = -4017.900000000001
$ws.Range("N136").Value = -17151

$ws = $wb.Worksheets.Item("CUL")
# Row 124: update H124, I124, J124, K124, L124, M124, N124
$ws.Range("H124").Value = 10910.625
$ws.Range("I124").Value = 1932.3334
$ws.Range("J124").Value = 16297.6
$ws.Range("K124").Value = 5797.0002
$ws.Range("L124").Value = 48892.8
$ws.Range("M124").Value = -887.0002000000004
$ws.Range("N124").Value = -58712.8
# Row 132: update H132, J132, L132, N132
$ws.Range("H132").Value = 1706.238
$ws.Range("J132").Value = 1943.1538
$ws.Range("L132").Value = 17488.3842
$ws.Range("N132").Value = -22548.3842

$ws = $wb.Worksheets.Item("GSM")
# Row 102: update H102, I102, J102, K102, L102, M102, N102
$ws.Range("H102").Value = 2322.2307
$ws.Range("I102").Value = 1868.9
$ws.Range("J102").Value = 3833.3333
$ws.Range("K102").Value = 1868.9
$ws.Range("L102").Value = 3833.3333
$ws.Range("M102").Value = -246.9000000000001
$ws.Range("N102").Value = -7077.3333
# Row 132: update H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 3800.6
$ws.Range("I132").Value = 3010.0417
$ws.Range("J132").Value = 4986.4375
$ws.Range("K132").Value = 9030.125100000001
$ws.Range("L132").Value = 14959.3125
$ws.Range("M132").Value = -6500.125100000001
$ws.Range("N132").Value = -20019.3125

$ws = $wb.Worksheets.Item("LTW")
# Row 132: update H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 2842.9443
$ws.Range("I132").Value = 2531.5
$ws.Range("J132").Value = 3465.8333
$ws.Range("K132").Value = 7594.5
$ws.Range("L132").Value = 10397.4999
$ws.Range("M132").Value = -5064.5
$ws.Range("N132").Value = -15457.4999
# Row 136: update H136, J136, L136, N136
$ws.Range("H136").Value = 3586.2856
$ws.Range("J136").Value = 4846.9
$ws.Range("L136").Value = 14540.7
$ws.Range("N136").Value = -19640.7

$ws = $wb.Worksheets.Item("WVR")
# Row 62: update H62, I62, K62, M62
$ws.Range("H62").Value = 4905.231
$ws.Range("I62").Value = 4121
$ws.Range("K62").Value = 4121
$ws.Range("M62").Value = -3497
# Row 65: update H65, I65, K65, M65
$ws.Range("H65").Value = 4905.231
$ws.Range("I65").Value = 4121
$ws.Range("K65").Value = 20605
$ws.Range("M65").Value = -17485
# Row 132: update H132, I132, K132, M132
$ws.Range("H132").Value = 2689.1462
$ws.Range("I132").Value = 2179.3142
$ws.Range("K132").Value = 6537.942599999999
$ws.Range("M132").Value = -4007.942599999999
# Row 136: update H136, J136, L136, N136
$ws.Range("H136").Value = 2579.3845
$ws.Range("J136").Value = 3892.818
$ws.Range("L136").Value = 11678.454
$ws.Range("N136").Value = -16778.454
